# Update Excel with new issue #14
# Appends a new row to the "Issues" sheet with the issue's data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

$row = 10

# Column A holds a numeric-looking issue id ("14") that must be stored as
# text (matching the rest of the sheet). Force text formatting while
# writing it, then clear the format again so no stray per-cell style is
# left behind - only the cell's text value should change.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "14"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "dfdddg. newrelic_alert_channel dfdfdfdf"
$ws.Cells.Item($row, 3).Value = "open"
$ws.Cells.Item($row, 4).Value = "2025-03-24T09:11:58Z"
$ws.Cells.Item($row, 5).Value = "enhancement"
